# API_Methoden.xlsx - add FlugzeugHandler methods (getAllRelationen, createFlugzeug,
# getAllFlugzeuge, assignFlugzeugToFlug) to the method-documentation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell values -------------------------------------------------
# NOTE: values are written in the exact order needed so that the
# resulting shared-string table indices line up with the target file
# (Excel appends a new shared string the first time a new text value is
# used, so order of assignment = order of new <si> entries).

# 1) RelationHandler.getAllRelationen (row 10) + "nicht implementiert" marker on
#    FlugHandler.createFlug (row 9, column G) + the getAllRelationen comment.
$ws.Range("B10").Value = "getAllRelationen"
$ws.Range("G9").Value = "nicht implementiert"
$ws.Range("G10").Value = 'gibt alle Relationen aus. Beispiel:  "1. Relation: Startort: FRA, Zielort: BOM (1500 km, 10:30:00 Stunden)"'

# 2) FlugzeugHandler.createFlugzeug (row 12)
$ws.Range("A12").Value = "FlugzeugHandler"
$ws.Range("B12").Value = "createFlugzeug"
$ws.Range("C12").Value = "String hersteller, String typ, int sitzplaetze"
$ws.Range("F12").Value = '"Erfolgreiche Anlage des Flugzeugs!"'
$ws.Range("E12").Value = '"Bitte geben Sie alle Information an!"'
$ws.Range("D12").Value = "String"

# 3) FlugzeugHandler.getAllFlugzeuge (row 13)
$ws.Range("B13").Value = "getAllFlugzeuge"
$ws.Range("C13").Value = "-"
$ws.Range("D13").Value = "List<String>"
$ws.Range("E13").Value = "-"
$ws.Range("F13").Value = "-"
$ws.Range("G13").Value = 'gibt alle Flugzeuge aus. Beispiel: "1. Flugzeug: Airbus A380-800 (853 Sitzplätze)"'

# 4) FlugzeugHandler.assignFlugzeugToFlug (row 14) - not implemented yet
$ws.Range("B14").Value = "assignFlugzeugToFlug"
$ws.Range("G14").Value = "nicht implementiert"

# Remaining reused values for row 10 (existing shared strings "-" / "List<String>")
$ws.Range("C10").Value = "-"
$ws.Range("D10").Value = "List<String>"
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "-"

# --- Formatting --------------------------------------------------------

# "nicht implementiert" cells get a red font, left/top aligned.
$ws.Range("G9").Font.Color = 255
$ws.Range("G14").Font.Color = 255

# Row 11 becomes a spacer row styled like the other "Schlecht" spacer rows
# (rows 5 and 8) - copy the formatting from row 8 so the existing style is
# reused instead of a new one being minted.
$ws.Range("A8:G8").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column F widened to fit the new "nicht implementiert" / comment text and no
# longer auto-fitted.
$ws.Columns.Item(6).ColumnWidth = 35.5

# --- View state ----------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("C2").Select()
